$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 4363.8887
$ws.Range("I76").Value = 4490.857
$ws.Range("J76").Value = 3919.5
$ws.Range("K76").Value = 4490.857
$ws.Range("L76").Value = 3919.5
$ws.Range("M76").Value = -4175.857
$ws.Range("N76").Value = -4549.5
$ws.Range("H79").Value = 4363.8887
$ws.Range("I79").Value = 4490.857
$ws.Range("J79").Value = 3919.5
$ws.Range("K79").Value = 4490.857
$ws.Range("L79").Value = 3919.5
$ws.Range("M79").Value = -3398.857
$ws.Range("N79").Value = -6103.5
$ws.Range("H86").Value = 4887.385
$ws.Range("I86").Value = 4287.2856
$ws.Range("K86").Value = 4287.2856
$ws.Range("M86").Value = -3164.2856
$ws.Range("H87").Value = 19999.834
$ws.Range("J87").Value = 19999.834
$ws.Range("L87").Value = 19999.834
$ws.Range("N87").Value = -22495.834
$ws.Range("H89").Value = 4887.385
$ws.Range("I89").Value = 4287.2856
$ws.Range("K89").Value = 21436.428
$ws.Range("M89").Value = -15820.428
$ws.Range("H90").Value = 19999.834
$ws.Range("J90").Value = 19999.834
$ws.Range("L90").Value = 59999.50199999999
$ws.Range("N90").Value = -72479.50199999999
$ws.Range("H101").Value = 1456.4
$ws.Range("I101").Value = 534.8570999999999
$ws.Range("K101").Value = 1604.5713
$ws.Range("M101").Value = 17.42870000000016
$ws.Range("H107").Value = 389.3
$ws.Range("J107").Value = 390.5
$ws.Range("L107").Value = 390.5
$ws.Range("N107").Value = -4230.5
$ws.Range("H137").Value = 2418.1
$ws.Range("J137").Value = 4124
$ws.Range("L137").Value = 12372
$ws.Range("N137").Value = -17472
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 202.38461
$ws.Range("I4").Value = 235.9
$ws.Range("K4").Value = 235.9
$ws.Range("M4").Value = -119.9
$ws.Range("H32").Value = 12258.118
$ws.Range("I32").Value = 9739.578
$ws.Range("J32").Value = 20353.428
$ws.Range("K32").Value = 9739.578
$ws.Range("L32").Value = 20353.428
$ws.Range("M32").Value = -9452.578
$ws.Range("N32").Value = -20927.428
$ws.Range("H45").Value = 6167.091
$ws.Range("I45").Value = 4764.8
$ws.Range("J45").Value = 7335.6665
$ws.Range("K45").Value = 4764.8
$ws.Range("L45").Value = 7335.6665
$ws.Range("M45").Value = -4387.8
$ws.Range("N45").Value = -8089.6665
$ws.Range("H61").Value = 4271.407
$ws.Range("I61").Value = 3048.5715
$ws.Range("J61").Value = 5588.3076
$ws.Range("K61").Value = 3048.5715
$ws.Range("L61").Value = 5588.3076
$ws.Range("M61").Value = -2836.5715
$ws.Range("N61").Value = -6012.3076
$ws.Range("H74").Value = 2156.8928
$ws.Range("I74").Value = 1894.6
$ws.Range("K74").Value = 1894.6
$ws.Range("M74").Value = -1020.6
$ws.Range("H77").Value = 2156.8928
$ws.Range("I77").Value = 1894.6
$ws.Range("K77").Value = 9473
$ws.Range("M77").Value = -5105
$ws.Range("H88").Value = 2083.842
$ws.Range("I88").Value = 2340.4443
$ws.Range("J88").Value = 1852.9
$ws.Range("K88").Value = 2340.4443
$ws.Range("L88").Value = 1852.9
$ws.Range("M88").Value = -1934.4443
$ws.Range("N88").Value = -2664.9
$ws.Range("H91").Value = 2083.842
$ws.Range("I91").Value = 2340.4443
$ws.Range("J91").Value = 1852.9
$ws.Range("K91").Value = 2340.4443
$ws.Range("L91").Value = 1852.9
$ws.Range("M91").Value = -936.4443000000001
$ws.Range("N91").Value = -4660.9
$ws.Range("H95").Value = 66666.664
$ws.Range("J95").Value = 66666.664
$ws.Range("L95").Value = 66666.664
$ws.Range("N95").Value = -72158.664
$ws.Range("H132").Value = 35718804
$ws.Range("I132").Value = 52635476
$ws.Range("K132").Value = 157906428
$ws.Range("M132").Value = -157903898
$ws.Range("H136").Value = 4271.407
$ws.Range("I136").Value = 3048.5715
$ws.Range("J136").Value = 5588.3076
$ws.Range("K136").Value = 9145.7145
$ws.Range("L136").Value = 16764.9228
$ws.Range("M136").Value = -6595.7145
$ws.Range("N136").Value = -21864.9228
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2947
$ws.Range("I105").Value = 1791.3636
$ws.Range("K105").Value = 1791.3636
$ws.Range("M105").Value = -44.36359999999991
$ws.Range("H134").Value = 4005.7334
$ws.Range("I134").Value = 2607.5789
$ws.Range("J134").Value = 6420.727
$ws.Range("K134").Value = 7822.736699999999
$ws.Range("L134").Value = 19262.181
$ws.Range("M134").Value = -5287.736699999999
$ws.Range("N134").Value = -24332.181
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4247.516
$ws.Range("I31").Value = 1025.1666
$ws.Range("K31").Value = 1025.1666
$ws.Range("M31").Value = -730.1666
$ws.Range("H34").Value = 4247.516
$ws.Range("I34").Value = 1025.1666
$ws.Range("K34").Value = 1025.1666
$ws.Range("M34").Value = -823.1666
$ws.Range("H122").Value = 37039580
$ws.Range("J122").Value = 4139.75
$ws.Range("L122").Value = 12419.25
$ws.Range("N122").Value = -17319.25
$ws.Range("H132").Value = 4911.7144
$ws.Range("I132").Value = 3376.4
$ws.Range("K132").Value = 10129.2
$ws.Range("M132").Value = -7599.200000000001
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 733.3333
$ws.Range("I132").Value = 600
$ws.Range("J132").Value = 800
$ws.Range("K132").Value = 5400
$ws.Range("L132").Value = 7200
$ws.Range("M132").Value = -2870
$ws.Range("N132").Value = -12260
$ws.Range("H137").Value = 2127.923
$ws.Range("J137").Value = 3474.75
$ws.Range("L137").Value = 10424.25
$ws.Range("N137").Value = -20624.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 164727.58
$ws.Range("I70").Value = 225818.6
$ws.Range("J70").Value = 12000
$ws.Range("K70").Value = 225818.6
$ws.Range("L70").Value = 12000
$ws.Range("M70").Value = -225548.6
$ws.Range("N70").Value = -12540
$ws.Range("H73").Value = 164727.58
$ws.Range("I73").Value = 225818.6
$ws.Range("J73").Value = 12000
$ws.Range("K73").Value = 225818.6
$ws.Range("L73").Value = 12000
$ws.Range("M73").Value = -224882.6
$ws.Range("N73").Value = -13872
$ws.Range("H102").Value = 24847246
$ws.Range("I102").Value = 3761292
$ws.Range("J102").Value = 125005530
$ws.Range("K102").Value = 3761292
$ws.Range("L102").Value = 125005530
$ws.Range("M102").Value = -3759670
$ws.Range("N102").Value = -125008774
$ws.Range("H132").Value = 4267.1724
$ws.Range("I132").Value = 3274
$ws.Range("K132").Value = 9822
$ws.Range("M132").Value = -7292
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 4206.6665
$ws.Range("I82").Value = 2160.818
$ws.Range("J82").Value = 9832.75
$ws.Range("K82").Value = 2160.818
$ws.Range("L82").Value = 9832.75
$ws.Range("M82").Value = -1799.818
$ws.Range("N82").Value = -10554.75
$ws.Range("H85").Value = 4206.6665
$ws.Range("I85").Value = 2160.818
$ws.Range("J85").Value = 9832.75
$ws.Range("K85").Value = 2160.818
$ws.Range("L85").Value = 9832.75
$ws.Range("M85").Value = -912.8180000000002
$ws.Range("N85").Value = -12328.75
$ws.Range("H110").Value = 38661.332
$ws.Range("J110").Value = 38661.332
$ws.Range("L110").Value = 38661.332
$ws.Range("N110").Value = -46841.332
$ws.Range("H136").Value = 4515.95
$ws.Range("I136").Value = 3263.52
$ws.Range("K136").Value = 9790.559999999999
$ws.Range("M136").Value = -7240.559999999999
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1968.3019
$ws.Range("I132").Value = 1410.0952
$ws.Range("J132").Value = 4099.636
$ws.Range("K132").Value = 4230.2856
$ws.Range("L132").Value = 12298.908
$ws.Range("M132").Value = -1700.2856
$ws.Range("N132").Value = -17358.908
